# Update "From Date" (C) and "End Date" (D) columns on Sheet1 to reflect
# the new reporting period (shifted forward by 5 months: Dec/Jan -> Feb/Mar 2026).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = 46082
$ws.Range("D3").Value = 46234

# Row 4
$ws.Range("C4").Value = 46054
$ws.Range("D4").Value = 46234

# Row 5
$ws.Range("C5").Value = 46054
$ws.Range("D5").Value = 46234

# Row 6
$ws.Range("C6").Value = 46082
$ws.Range("D6").Value = 46234

# Row 7
$ws.Range("C7").Value = 46082
$ws.Range("D7").Value = 46234

# Row 8
$ws.Range("C8").Value = 46082
$ws.Range("D8").Value = 46173

# Row 9
$ws.Range("C9").Value = 46082
$ws.Range("D9").Value = 46112

# Row 10
$ws.Range("C10").Value = 46082
$ws.Range("D10").Value = 46112

# Row 11
$ws.Range("C11").Value = 46082
$ws.Range("D11").Value = 46173

# Row 12
$ws.Range("C12").Value = 46082
$ws.Range("D12").Value = 46173

# Row 13
$ws.Range("C13").Value = 46082
$ws.Range("D13").Value = 46173

# Row 14
$ws.Range("C14").Value = 46082
$ws.Range("D14").Value = 46112

# Row 15
$ws.Range("C15").Value = 46082
$ws.Range("D15").Value = 46112
